$wb = $excel.ActiveWorkbook

# Values to update in column F for rows 2-6 and 8 (想去人数 - number of people interested)
$updates = @{
    2 = 2245
    3 = 1702
    4 = 333
    5 = 1085
    6 = 797
    8 = 5821
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
